$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Re-ordena el detalle de periodos de mora por trabajador/periodo (orden cronologico)
$data = @(
    ,@("CC", "73230194", "ABEL FERNANDO BRIEVA RODRIGUEZ", "1608", 27578)
    ,@("CC", "73230194", "ABEL FERNANDO BRIEVA RODRIGUEZ", "1609", 27578)
    ,@("CC", "73230194", "ABEL FERNANDO BRIEVA RODRIGUEZ", "1610", 27578)
    ,@("CC", "1052988508", "DUGLAS ROBLES NIÑO", "1805", 1042)
    ,@("CC", "1052988508", "DUGLAS ROBLES NIÑO", "1806", 31249)
    ,@("CC", "1052988508", "DUGLAS ROBLES NIÑO", "1807", 31249)
    ,@("CC", "1052988508", "DUGLAS ROBLES NIÑO", "1808", 31249)
    ,@("CC", "1052988508", "DUGLAS ROBLES NIÑO", "1809", 31249)
    ,@("CC", "1052988508", "DUGLAS ROBLES NIÑO", "1810", 31249)
    ,@("CC", "73230194", "ABEL FERNANDO BRIEVA RODRIGUEZ", "1811", 31249)
    ,@("CC", "1047445886", "SUSANA KATHERINE ROJAS DE LA ROSA", "1811", 31249)
    ,@("CC", "1052988508", "DUGLAS ROBLES NIÑO", "1811", 31249)
    ,@("CC", "73230194", "ABEL FERNANDO BRIEVA RODRIGUEZ", "1812", 31249)
    ,@("CC", "1047445886", "SUSANA KATHERINE ROJAS DE LA ROSA", "1812", 31249)
    ,@("CC", "1052988508", "DUGLAS ROBLES NIÑO", "1812", 31249)
    ,@("CC", "73230194", "ABEL FERNANDO BRIEVA RODRIGUEZ", "1901", 31249)
    ,@("CC", "1047445886", "SUSANA KATHERINE ROJAS DE LA ROSA", "1901", 31249)
    ,@("CC", "1052988508", "DUGLAS ROBLES NIÑO", "1901", 31249)
    ,@("CC", "73230194", "ABEL FERNANDO BRIEVA RODRIGUEZ", "1902", 31249)
    ,@("CC", "1047445886", "SUSANA KATHERINE ROJAS DE LA ROSA", "1902", 31249)
    ,@("CC", "1052988508", "DUGLAS ROBLES NIÑO", "1902", 31249)
    ,@("CC", "73230194", "ABEL FERNANDO BRIEVA RODRIGUEZ", "1903", 26041)
    ,@("CC", "1047445886", "SUSANA KATHERINE ROJAS DE LA ROSA", "1903", 26041)
    ,@("CC", "1052988508", "DUGLAS ROBLES NIÑO", "1903", 26041)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 16 + $i
    $ws.Range("B$row").Value = $data[$i][0]
    $ws.Range("C$row").Value = $data[$i][1]
    $ws.Range("D$row").Value = $data[$i][2]
    $ws.Range("E$row").Value = $data[$i][3]
    $ws.Range("F$row").Value = $data[$i][4]
}
